$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z8").Value = "2025-10-17T07:09:39.336865"
$ws.Range("Z9:Z19").Value = "2025-10-17T07:09:39.337866"
$ws.Range("Z20:Z23").Value = "2025-10-17T07:09:39.338915"
$ws.Range("Z24:Z30").Value = "2025-10-17T07:09:39.339423"
$ws.Range("Z31:Z39").Value = "2025-10-17T07:09:39.340669"
$ws.Range("Z40:Z45").Value = "2025-10-17T07:09:39.341670"
$ws.Range("Z46:Z74").Value = "2025-10-17T07:09:39.429492"
$ws.Range("Z75").Value = "2025-10-17T07:09:39.523170"
$ws.Range("Z76:Z78").Value = "2025-10-17T07:09:39.524088"
$ws.Range("Z79:Z82").Value = "2025-10-17T07:09:39.525091"
$ws.Range("Z83").Value = "2025-10-17T07:09:39.526087"
$ws.Range("Z84").Value = "2025-10-17T07:09:39.527088"
$ws.Range("Z85:Z87").Value = "2025-10-17T07:09:39.527354"
$ws.Range("Z88:Z92").Value = "2025-10-17T07:09:39.527895"
$ws.Range("Z93:Z98").Value = "2025-10-17T07:09:39.528886"
$ws.Range("Z99").Value = "2025-10-17T07:09:39.529887"
$ws.Range("Z100:Z102").Value = "2025-10-17T07:09:39.530022"
$ws.Range("Z103:Z107").Value = "2025-10-17T07:09:39.607568"
$ws.Range("Z108:Z112").Value = "2025-10-17T07:09:39.624096"
